$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 with new values
$ws.Range("A2").Value = 45033.50694444445
$ws.Range("B2").Value = 5.378
$ws.Range("C2").Value = 5.209
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 13.173
$ws.Range("F2").Value = 11.188
$ws.Range("G2").Value = 5.041
$ws.Range("H2").Value = 13.381
$ws.Range("I2").Value = 7.973
$ws.Range("J2").Value = 3.774
$ws.Range("K2").Value = 6.167
$ws.Range("L2").Value = 6.371
$ws.Range("M2").Value = 6.162
$ws.Range("N2").Value = 1.68
$ws.Range("O2").Value = 4.805
$ws.Range("P2").Value = 6.849
$ws.Range("Q2").Value = 4.447
$ws.Range("R2").Value = 0.297
$ws.Range("S2").Value = 0.425
$ws.Range("T2").Value = 72.19199999999999
$ws.Range("U2").Value = 14.308
$ws.Range("V2").Value = 4.642
$ws.Range("W2").Value = 8.505000000000001
$ws.Range("X2").Value = 6.25
$ws.Range("Y2").Value = 0.6899999999999999
$ws.Range("Z2").Value = 6.977
$ws.Range("AA2").Value = 4.068
$ws.Range("AB2").Value = 5.066
$ws.Range("AC2").Value = 6.321
$ws.Range("AD2").Value = 6.735
$ws.Range("AE2").Value = 0.773
$ws.Range("AF2").Value = 11.064
$ws.Range("AG2").Value = 3.76
$ws.Range("AH2").Value = 5.189
$ws.Range("A3").Value = 45033.51388888889
$ws.Range("B3").Value = 19.9
$ws.Range("C3").Value = 15.468
$ws.Range("D3").Value = 0.466
$ws.Range("E3").Value = 44.248
$ws.Range("F3").Value = 36.535
$ws.Range("G3").Value = 16.067
$ws.Range("H3").Value = 57.758
$ws.Range("I3").Value = 25.022
$ws.Range("J3").Value = 11.45
$ws.Range("K3").Value = 17.254
$ws.Range("L3").Value = 18.385
$ws.Range("M3").Value = 18.987
$ws.Range("N3").Value = 5.192
$ws.Range("O3").Value = 15.826
$ws.Range("P3").Value = 22.88
$ws.Range("Q3").Value = 13.431
$ws.Range("R3").Value = 0.283
$ws.Range("S3").Value = 0.754
$ws.Range("T3").Value = 238.458
$ws.Range("U3").Value = 45.119
$ws.Range("V3").Value = 14.759
$ws.Range("W3").Value = 29.971
$ws.Range("X3").Value = 16.66
$ws.Range("Y3").Value = 2.119
$ws.Range("Z3").Value = 29.232
$ws.Range("AA3").Value = 13.107
$ws.Range("AB3").Value = 12.323
$ws.Range("AC3").Value = 14.551
$ws.Range("AD3").Value = 19.337
$ws.Range("AE3").Value = 0.419
$ws.Range("AF3").Value = 52.403
$ws.Range("AG3").Value = 8.923
$ws.Range("AH3").Value = 18.276
$ws.Range("A4").Value = 45033.52083333334
$ws.Range("B4").Value = 2.344
$ws.Range("C4").Value = 2.088
$ws.Range("D4").Value = 0.005
$ws.Range("E4").Value = 5.787
$ws.Range("F4").Value = 4.712
$ws.Range("G4").Value = 2.122
$ws.Range("H4").Value = 15.662
$ws.Range("I4").Value = 3.536
$ws.Range("J4").Value = 1.847
$ws.Range("K4").Value = 2.696
$ws.Range("L4").Value = 2.787
$ws.Range("M4").Value = 2.571
$ws.Range("N4").Value = 0.745
$ws.Range("O4").Value = 1.985
$ws.Range("P4").Value = 3.254
$ws.Range("Q4").Value = 1.891
$ws.Range("R4").Value = 0.101
$ws.Range("S4").Value = 0.187
$ws.Range("T4").Value = 26.914
$ws.Range("U4").Value = 6.761
$ws.Range("V4").Value = 1.951
$ws.Range("W4").Value = 4.292
$ws.Range("X4").Value = 2.746
$ws.Range("Y4").Value = 0.296
$ws.Range("Z4").Value = 7.187
$ws.Range("AA4").Value = 1.803
$ws.Range("AB4").Value = 2.136
$ws.Range("AC4").Value = 2.507
$ws.Range("AD4").Value = 2.932
$ws.Range("AE4").Value = 0.281
$ws.Range("AF4").Value = 14.831
$ws.Range("AG4").Value = 1.444
$ws.Range("AH4").Value = 2.367
$ws.Range("A5").Value = 45033.52777777778
$ws.Range("B5").Value = 11.59
$ws.Range("C5").Value = 8.94
$ws.Range("D5").Value = 0.31
$ws.Range("E5").Value = 25.77
$ws.Range("F5").Value = 21.22
$ws.Range("G5").Value = 9.33
$ws.Range("H5").Value = 34.17
$ws.Range("I5").Value = 14.6
$ws.Range("J5").Value = 6.65
$ws.Range("K5").Value = 10.04
$ws.Range("L5").Value = 10.7
$ws.Range("M5").Value = 11.01
$ws.Range("N5").Value = 3.02
$ws.Range("O5").Value = 9.17
$ws.Range("P5").Value = 13.33
$ws.Range("Q5").Value = 7.8
$ws.Range("R5").Value = 0.16
$ws.Range("S5").Value = 0.43
$ws.Range("T5").Value = 135.55
$ws.Range("U5").Value = 26.15
$ws.Range("V5").Value = 8.57
$ws.Range("W5").Value = 17.39
$ws.Range("X5").Value = 9.640000000000001
$ws.Range("Y5").Value = 1.23
$ws.Range("Z5").Value = 16.77
$ws.Range("AA5").Value = 7.64
$ws.Range("AB5").Value = 7.13
$ws.Range("AC5").Value = 8.380000000000001
$ws.Range("AD5").Value = 11.26
$ws.Range("AE5").Value = 0.21
$ws.Range("AF5").Value = 30.64
$ws.Range("AG5").Value = 5.13
$ws.Range("AH5").Value = 10.67

# Remove row 6 entirely (dataset shrunk from 5 to 4 data rows)
$ws.Rows.Item(6).Delete()

# Adjust column widths
$ws.Columns.Item(2).ColumnWidth = 6.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 6.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667
